$d = $word.ActiveDocument

# Locate the "Practice book 16- Test 2" heading, then the first
# "Band Score:" paragraph that follows it (the Writing Task 1 score
# line for that test, which currently has no score typed after it).
$heading = $d.Content
$headingFound = $heading.Find.Execute("Practice book 16- Test 2", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $headingFound) {
    throw "Could not find 'Practice book 16- Test 2' heading"
}

$scoreRange = $d.Range($heading.End, $d.Content.End)
$scoreFound = $scoreRange.Find.Execute("Band Score:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $scoreFound) {
    throw "Could not find 'Band Score:' paragraph after the heading"
}

# Append the band score right after the existing "Band Score:" text.
# The inserted run inherits the same bold, size-24 formatting already
# used for that "Band Score:" label.
$scoreRange.Collapse(0)
$scoreRange.InsertAfter(" 6.5 – 7")
